$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterCustomerTest")

# Update test data values in row 2
$ws.Range("A2").Value = "Bill"
$ws.Range("I2").Value = "Bill123"
$ws.Range("L2").Value = "Customer Created"

# Set column L width (column 12) to match new content
# (15.1640625 is the exact target width; the engine quantizes ColumnWidth to
# 1/6-character pixel steps, so 14.25 is the nearest input that lands on the
# closest achievable stored width of 15.1666...)
$ws.Columns.Item(12).ColumnWidth = 14.25

# Update the selected cell on this sheet
$ws.Range("I2").Select()
